# Update the DEQM Reporter Client Capability Statement workbook from the
# STU3 / 1.0.0 era to the R4 / 1.1.0 era ("update to r4 qa").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "meta" sheet: bump version / fhirVersion / ig URL
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("meta")

$meta.Range("B3").Value = "1.1.0"   # version
$meta.Range("B4").Value = "4.0.0"   # fhirVersion (STU3 "3.0.1" -> R4 "4.0.0")
$meta.Range("B6").Value = "http://hl7.org/fhir/us/davinci-deqm/ImplementationGuide/hl7.fhir.us.davinci-deqm-1.1.0"   # ig

# ---------------------------------------------------------------------
# 2. "profiles" sheet: drop the "/STU3" path segment from every profile URL
# ---------------------------------------------------------------------
$profiles = $wb.Worksheets.Item("profiles")

$profileUrls = @(
    "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/devicerequest-deqm",
    "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/medicationadministration-deqm",
    "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/summary-measurereport-deqm",
    "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/practitioner-deqm",
    "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/deviceusestatement-deqm",
    "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/organization-deqm",
    "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/coverage-deqm",
    "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/indv-measurereport-deqm",
    "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/medicationrequest-deqm"
)

for ($i = 0; $i -lt $profileUrls.Length; $i++) {
    $row = $i + 2
    $profiles.Cells.Item($row, 1).Value = $profileUrls[$i]
}

# ---------------------------------------------------------------------
# 3. Selections / active sheet bookkeeping, matching the saved state in
#    the workbook (meta!B6 selected, ops tab no longer active, profiles
#    becomes the active tab with A19 selected).
# ---------------------------------------------------------------------
$meta.Activate() | Out-Null
$meta.Range("B6").Select() | Out-Null

$ops = $wb.Worksheets.Item("ops")
$ops.Activate() | Out-Null
$ops.Range("F5").Select() | Out-Null

$profiles.Activate() | Out-Null
$profiles.Range("A19").Select() | Out-Null

Write-Host "Updated meta/profiles sheets to R4 1.1.0 references"
